# Give the slide master / layout titles (and the shapes below them) more
# vertical room - especially useful on the chart layout. Point values below
# are the EMU targets (839788, 365125, 10515600, 974000, ...) converted to
# points (EMU / 12700); a couple of them are nudged by a hair so that the
# host's float32 round-trip lands back on the exact EMU the diff expects.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# ---------------------------------------------------------------------
# CustomLayout 5 = "Vergelijking" (Comparison)
# ---------------------------------------------------------------------
$lay5 = $m.CustomLayouts.Item(5)

# Title 1: ext cy 774000 -> 898410
$sh = $lay5.Shapes.Item(1)
$sh.Left   = 66.12504577636719
$sh.Top    = 28.75
$sh.Width  = 828
$sh.Height = 70.74095153808594

# Text Placeholder 2: off 839788,1185863 -> 839788,1263535 ; ext 5157787,823912 -> 5157787,746240
$sh = $lay5.Shapes.Item(2)
$sh.Left   = 66.12504577636719
$sh.Top    = 99.49095153808594
$sh.Width  = 406.12496062992125
$sh.Height = 58.75905511811023

# Text Placeholder 4: off 6172200,1185863 -> 6172200,1263535 ; ext 5183188,823912 -> 5183188,746240
$sh = $lay5.Shapes.Item(4)
$sh.Left   = 486.0
$sh.Top    = 99.49095153808594
$sh.Width  = 408.12506103515625
$sh.Height = 58.75905511811023

# ---------------------------------------------------------------------
# CustomLayout 15 = "1_Vergelijking"
# ---------------------------------------------------------------------
$lay15 = $m.CustomLayouts.Item(15)

# Title 1: off 839788,365125 -> 839788,365124 ; ext 10515600,774000 -> 10515600,823911
$sh = $lay15.Shapes.Item(1)
$sh.Left   = 66.12504577636719
$sh.Top    = 28.74992125984252
$sh.Width  = 828
$sh.Height = 64.87488188976378

# ---------------------------------------------------------------------
# CustomLayout 16 = "2_Vergelijking"
# ---------------------------------------------------------------------
$lay16 = $m.CustomLayouts.Item(16)

# Title 1: ext 10515600,774000 -> 10515600,823912 (offset unchanged)
$sh = $lay16.Shapes.Item(1)
$sh.Left   = 66.12504577636719
$sh.Top    = 28.75
$sh.Width  = 828
$sh.Height = 64.87496062992126

# ---------------------------------------------------------------------
# CustomLayout 17 = "3_Titel en object" (chart slide)
# ---------------------------------------------------------------------
$lay17 = $m.CustomLayouts.Item(17)

# Title 1: previously inherited (<p:spPr/>), now explicit off 838200,365125 ext 10515600,973223
$sh = $lay17.Shapes.Item(1)
$sh.Left   = 66.0
$sh.Top    = 28.75
$sh.Width  = 828
$sh.Height = 76.6317367553711

# Chart placeholder: off 838200,1226012 -> 838200,1413163 ; ext 10515600,5470525 -> 10515600,5283373
$sh = $lay17.Shapes.Item(2)
$sh.Left   = 66.0
$sh.Top    = 111.2726821899414
$sh.Width  = 828
$sh.Height = 416.0136413574219

# ---------------------------------------------------------------------
# Slide Master itself
# ---------------------------------------------------------------------

# Title Placeholder 1: ext 10515600,773863 -> 10515600,902200 (offset unchanged)
$sh = $m.Shapes.Item(1)
$sh.Left   = 66.0
$sh.Top    = 28.750080108642578
$sh.Width  = 828
$sh.Height = 71.03937530517578

# Text Placeholder 2: off 838200,1267326 -> 838200,1354974 ; ext 10515600,5499234 -> 10515600,5411585
$sh = $m.Shapes.Item(2)
$sh.Left   = 66.0
$sh.Top    = 106.69087219238281
$sh.Width  = 828
$sh.Height = 426.10905511811023
